# Populate the "Catcher's View" pitch-by-pitch table (Pitch / Choice / Result
# columns F/G/H) and a few derived Outcome cells (M column) that were left
# blank in the template, plus normalize the "Pitch Mix" (J column) ordering.
# This builds the data needed for the hitter's strikezone visual.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- At-bat 1 (rows 10-14), pitch mix in J17 --------------------------------
$ws.Range("F10").Value = "FB"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Ball"
$ws.Range("M10").Value = "71.6 MPH"

$ws.Range("F11").Value = "FB"
$ws.Range("G11").Value = "Swing"
$ws.Range("H11").Value = "Foul"

$ws.Range("F12").Value = "CH"
$ws.Range("G12").Value = "Swing"
$ws.Range("H12").Value = "Strike"
$ws.Range("M12").Value = "22.4°"

$ws.Range("F13").Value = "CH"
$ws.Range("G13").Value = "Take"
$ws.Range("H13").Value = "Ball"

$ws.Range("F14").Value = "FB"
$ws.Range("G14").Value = "Swing"
$ws.Range("H14").Value = "In Play"

$ws.Range("J17").Value = "CH,CB,FB"

# ---- At-bat 2 (rows 19-25), pitch mix in J26 --------------------------------
$ws.Range("F19").Value = "FB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Strike"

$ws.Range("F20").Value = "FB"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Strike"

$ws.Range("F21").Value = "CH"
$ws.Range("G21").Value = "Take"
$ws.Range("H21").Value = "Ball"
$ws.Range("M21").Value = $null

$ws.Range("F22").Value = "CB"
$ws.Range("G22").Value = "Take"
$ws.Range("H22").Value = "Ball"

$ws.Range("F23").Value = "FB"
$ws.Range("G23").Value = "Swing"
$ws.Range("H23").Value = "Foul"

$ws.Range("F24").Value = "CB"
$ws.Range("G24").Value = "Take"
$ws.Range("H24").Value = "Ball"
$ws.Range("M24").Value = "Strikeout"

$ws.Range("F25").Value = "FB"
$ws.Range("G25").Value = "Swing"
$ws.Range("H25").Value = "Strike"

$ws.Range("J26").Value = "CH,CB,FB"

# ---- At-bat 3 (rows 28-34), pitch mix in J35 --------------------------------
$ws.Range("F28").Value = "CB"
$ws.Range("G28").Value = "Take"
$ws.Range("H28").Value = "Ball"
$ws.Range("M28").Value = "99.63 MPH"

$ws.Range("F29").Value = "CH"
$ws.Range("G29").Value = "Swing"
$ws.Range("H29").Value = "In Play"

$ws.Range("M30").Value = "10.11°"

$ws.Range("J35").Value = "CH,CB,FB,SL"

# ---- At-bat 4 (rows 37-43), pitch mix in J44 --------------------------------
$ws.Range("F37").Value = "FB"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Ball"

$ws.Range("F38").Value = "FB"
$ws.Range("G38").Value = "Take"
$ws.Range("H38").Value = "Strike"

$ws.Range("F39").Value = "FB"
$ws.Range("G39").Value = "Swing"
$ws.Range("H39").Value = "Foul"
$ws.Range("M39").Value = $null

$ws.Range("F40").Value = "FB"
$ws.Range("G40").Value = "Take"
$ws.Range("H40").Value = "Ball"

$ws.Range("F41").Value = "FB"
$ws.Range("G41").Value = "Swing"
$ws.Range("H41").Value = "Foul"

$ws.Range("F42").Value = "CH"
$ws.Range("G42").Value = "Take"
$ws.Range("H42").Value = "Ball"
$ws.Range("M42").Value = "Hit By Pitch"

$ws.Range("F43").Value = "SL"
$ws.Range("G43").Value = "Take"
$ws.Range("H43").Value = "HBP"

$ws.Range("J44").Value = "CH,FB,SL"

# ---- At-bat 5 (rows 46-52), pitch mix in J53 --------------------------------
$ws.Range("F46").Value = "SL"
$ws.Range("G46").Value = "Take"
$ws.Range("H46").Value = "HBP"
$ws.Range("M46").Value = "20.44 MPH"

$ws.Range("M48").Value = "7.81°"

$ws.Range("M51").Value = "Hit By Pitch"

$ws.Range("J53").Value = "CH,FB,SL"

# ---- At-bat 6 (rows 61-64), pitch mix in J68 --------------------------------
$ws.Range("F61").Value = "CH"
$ws.Range("G61").Value = "Take"
$ws.Range("H61").Value = "Strike"
$ws.Range("M61").Value = "53.02 MPH"

$ws.Range("F62").Value = "CH"
$ws.Range("G62").Value = "Take"
$ws.Range("H62").Value = "Ball"

$ws.Range("F63").Value = "CH"
$ws.Range("G63").Value = "Take"
$ws.Range("H63").Value = "Ball"
$ws.Range("M63").Value = "31.83°"

$ws.Range("F64").Value = "CH"
$ws.Range("G64").Value = "Swing"
$ws.Range("H64").Value = "In Play"

$ws.Range("J68").Value = "CH,FB,SL"
